$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Add the new "optimization_parameters" sheet after the last existing sheet
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "optimization_parameters"

# A custom cell style (mirrors the "Normal 2" style pasted in from another
# workbook in the original edit) carrying the Arial 10pt black font used
# for all the data on this sheet.
$dataStyle = $wb.Styles.Add("Normal 2")
$dataStyle.Font.Size = 10
$dataStyle.Font.Name = "Arial"
$dataStyle.Font.Color = 0

$ws.Range("A1:E18").Style = "Normal 2"
$ws.Columns.Item(1).ColumnWidth = 20.625

# ------------------------------------------------------------------
# Populate the optimization_parameters data
# ------------------------------------------------------------------
$ws.Range("A1").Value = "optimization_parameter"
$ws.Range("B1").Value = "value"

$ws.Range("A2").Value = "alpha"
$ws.Range("B2").Value = 0.002
$ws.Range("B2").NumberFormat = "0.00E+00"

$ws.Range("A3").Value = "kk_max"
$ws.Range("B3").Value = 1

$ws.Range("A4").Value = "MaxIter"
$ws.Range("B4").Value = 100000000
$ws.Range("B4").NumberFormat = "0.00E+00"

$ws.Range("A5").Value = "TolFun"
$ws.Range("B5").Value = 0.000001
$ws.Range("B5").NumberFormat = "0.00E+00"

$ws.Range("A6").Value = "MaxFunEval"
$ws.Range("B6").Value = 100000000
$ws.Range("B6").NumberFormat = "0.00E+00"

$ws.Range("A7").Value = "TolX"
$ws.Range("B7").Value = 0.000001
$ws.Range("B7").NumberFormat = "0.00E+00"

$ws.Range("A8").Value = "production_function"
$ws.Range("B8").Value = "Sigmoid"

$ws.Range("A9").Value = "L_curve"
$ws.Range("B9").Value = 0

$ws.Range("A10").Value = "estimate_params"
$ws.Range("B10").Value = 1

$ws.Range("A11").Value = "make_graphs"
$ws.Range("B11").Value = 1

$ws.Range("A12").Value = "fix_P"
$ws.Range("B12").Value = 0

$ws.Range("A13").Value = "fix_b"
$ws.Range("B13").Value = 0

$ws.Range("A14").Value = "expression_timepoints"
$ws.Range("B14").Value = 15
$ws.Range("C14").Value = 30
$ws.Range("D14").Value = 60

$ws.Range("A15").Value = "Strain"
$ws.Range("B15").Value = "wt"
$ws.Range("C15").Value = "dgln3"

$ws.Range("A16").Value = "simulation_timepoints"
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 10
$ws.Range("E16").Value = 15

$ws.Range("A17").Value = "species"
$ws.Range("B17").Value = "Saccharomyces cerevisiae"

$ws.Range("A18").Value = "taxon_id"
$ws.Range("B18").Value = 559292

# ------------------------------------------------------------------
# View state: new sheet becomes the active/selected tab, scrolled &
# selected the way the original was left.
# ------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 100
$ws.Range("C22").Select() | Out-Null

# Iterative-calculation settings used by the optimization routine.
$excel.MaxIterations = 100
$excel.MaxChange = 0.0001
